$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column C, row 1
$ws.Range("C1").Value = "Xóa nợ"

# Append new entries to B2 (Phan Ngọc Hiếu's notes)
$ws.Range("B2").Value = "Sáng 10/03: Đi trễ`nChiều 10/03: Về sớm`n13/03: Chưa chuẩn bị bài`n15/03: Đau xin nghỉ`n16/03: Đau xin nghỉ"

# Append new entries to B4 (Hưng's notes)
$ws.Range("B4").Value = "10/03: Xin về sớm`n13/03: Đi trễ`n16/03: Đi trễ`n16/03: Xin về sơm"

# New row 6: Thiện
$ws.Range("A6").Value = "Thiện"
$ws.Range("B6").Value = "16/03: Đi trể`n17/03: Xin nghỉ 1 ngày"
$ws.Range("B6").WrapText = $true

# Adjust row heights to match the updated content
# (row 5 already has ht=45 inherited from before, so no need to touch it)
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 63.75
$ws.Rows.Item(6).RowHeight = 30

# Update the active cell selection
$ws.Range("B7").Select()
